# Scheduled refresh of market-price-derived profit columns (H:N) across all
# job sheets, sourced from the latest Universalis average-price snapshot.
# Leve Item ID / Leve Name / Leve Level columns (A:G) are untouched.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 40: Stuck in the Moment
$ws.Range("H40").Value = 3535.7273
$ws.Range("J40").Value = 3687.125
$ws.Range("L40").Value = 3687.125
$ws.Range("N40").Value = -4037.125

# Row 62: The Mustache Suits Him
$ws.Range("H62").Value = 41770.785
$ws.Range("I62").Value = 4950
$ws.Range("J62").Value = 47907.582
$ws.Range("K62").Value = 4950
$ws.Range("L62").Value = 47907.582
$ws.Range("M62").Value = -4326
$ws.Range("N62").Value = -49155.582

# Row 65: Forgery of Convenience (L)
$ws.Range("H65").Value = 41770.785
$ws.Range("I65").Value = 4950
$ws.Range("J65").Value = 47907.582
$ws.Range("K65").Value = 24750
$ws.Range("L65").Value = 239537.91
$ws.Range("M65").Value = -21630
$ws.Range("N65").Value = -245777.91

# Row 88: The Grave of Hemlock Groves
$ws.Range("H88").Value = 1397.0834
$ws.Range("I88").Value = 1422
$ws.Range("J88").Value = 1384.625
$ws.Range("K88").Value = 1422
$ws.Range("L88").Value = 1384.625
$ws.Range("M88").Value = -1016
$ws.Range("N88").Value = -2196.625

# Row 91: Dappling the Highlands (L)
$ws.Range("H91").Value = 1397.0834
$ws.Range("I91").Value = 1422
$ws.Range("J91").Value = 1384.625
$ws.Range("K91").Value = 1422
$ws.Range("L91").Value = 1384.625
$ws.Range("M91").Value = -18
$ws.Range("N91").Value = -4192.625

# Row 116: Growing Up
$ws.Range("H116").Value = 34385104
$ws.Range("I116").Value = 20931570
$ws.Range("J116").Value = 66673584
$ws.Range("K116").Value = 20931570
$ws.Range("L116").Value = 66673584
$ws.Range("M116").Value = -20928128
$ws.Range("N116").Value = -66680468

# Row 132: Fast-forwarding Flora
$ws.Range("H132").Value = 3217.9307
$ws.Range("I132").Value = 3127.1226
$ws.Range("K132").Value = 9381.3678
$ws.Range("M132").Value = -6851.3678

# Row 136: I Like Big Brush and I Cannot Lie
$ws.Range("H136").Value = 59999
$ws.Range("J136").Value = 59999
$ws.Range("L136").Value = 59999
$ws.Range("N136").Value = -70199

$ws = $wb.Worksheets.Item("ARM")
# Row 61: Dealing with the Tough Stuff
$ws.Range("H61").Value = 2858630
$ws.Range("I61").Value = 2501349.2
$ws.Range("J61").Value = 3335004.2
$ws.Range("K61").Value = 2501349.2
$ws.Range("L61").Value = 3335004.2
$ws.Range("M61").Value = -2501137.2
$ws.Range("N61").Value = -3335428.2

# Row 136: Metal with Mettle
$ws.Range("H136").Value = 2858630
$ws.Range("I136").Value = 2501349.2
$ws.Range("J136").Value = 3335004.2
$ws.Range("K136").Value = 7504047.600000001
$ws.Range("L136").Value = 10005012.6
$ws.Range("M136").Value = -7501497.600000001
$ws.Range("N136").Value = -10010112.6

$ws = $wb.Worksheets.Item("BSM")
# Row 22: Riveting Run
$ws.Range("H22").Value = 196.66667
$ws.Range("I22").Value = 196.66667
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 196.66667
$ws.Range("L22").Value = 0
$ws.Range("N22").Value = -23.66667000000001

# Row 86: Through Thick and Thin
$ws.Range("H86").Value = 2082.2856
$ws.Range("I86").Value = 1680.25
$ws.Range("J86").Value = 2618.3333
$ws.Range("K86").Value = 1680.25
$ws.Range("L86").Value = 2618.3333
$ws.Range("M86").Value = -557.25
$ws.Range("N86").Value = -4864.3333

# Row 89: Piercing Eyes Deserve Piercing Shafts (L)
$ws.Range("H89").Value = 2082.2856
$ws.Range("I89").Value = 1680.25
$ws.Range("J89").Value = 2618.3333
$ws.Range("K89").Value = 8401.25
$ws.Range("L89").Value = 13091.6665
$ws.Range("M89").Value = -2785.25
$ws.Range("N89").Value = -24323.6665

$ws = $wb.Worksheets.Item("CRP")
# Row 31: Wall Not Found
$ws.Range("H31").Value = 2623.66
$ws.Range("J31").Value = 2856.045
$ws.Range("L31").Value = 2856.045
$ws.Range("N31").Value = -3446.045

# Row 34: Armoires of the Rich and Famous
$ws.Range("H34").Value = 2623.66
$ws.Range("J34").Value = 2856.045
$ws.Range("L34").Value = 2856.045
$ws.Range("N34").Value = -3260.045

# Row 35: Storm of Swords
$ws.Range("H35").Value = 636.1667
$ws.Range("I35").Value = 636.1667
$ws.Range("K35").Value = 636.1667
$ws.Range("M35").Value = -342.1667

# Row 94: Beech, Please
$ws.Range("H94").Value = 3828.0667
$ws.Range("I94").Value = 3928
$ws.Range("K94").Value = 3928
$ws.Range("M94").Value = -3477

# Row 99: O Pine
$ws.Range("H99").Value = 2685
$ws.Range("I99").Value = 2628
$ws.Range("K99").Value = 2628
$ws.Range("M99").Value = -1130

# Row 105: Zelkova, My Love
$ws.Range("H105").Value = 1242
$ws.Range("I105").Value = 1242
$ws.Range("K105").Value = 1242
$ws.Range("M105").Value = 505

# Row 122: Timber of Tenkonto
$ws.Range("H122").Value = 3513.4546
$ws.Range("I122").Value = 2072.111
$ws.Range("K122").Value = 6216.333
$ws.Range("M122").Value = -3766.333

# Row 126: A Better Conductor
$ws.Range("H126").Value = 2685
$ws.Range("I126").Value = 2628
$ws.Range("K126").Value = 7884
$ws.Range("M126").Value = -5414

$ws = $wb.Worksheets.Item("CUL")
# Row 12: Butter Me Up
$ws.Range("H12").Value = 74
$ws.Range("J12").Value = 80.5
$ws.Range("L12").Value = 241.5
$ws.Range("N12").Value = -587.5

# Row 34: Fever Pitch
$ws.Range("H34").Value = 1769.75
$ws.Range("J34").Value = 2628
$ws.Range("L34").Value = 7884
$ws.Range("N34").Value = -8052

$ws = $wb.Worksheets.Item("GSM")
# Row 53: North Ore South
$ws.Range("H53").Value = 21663
$ws.Range("I53").Value = 21663
$ws.Range("K53").Value = 21663
$ws.Range("M53").Value = -21032

# Row 58: The Big Red
$ws.Range("H58").Value = 15516.667
$ws.Range("I58").Value = 11380
$ws.Range("K58").Value = 11380
$ws.Range("M58").Value = -11103

# Row 122: Awarding Academic Excellence
$ws.Range("H122").Value = 3415.1875
$ws.Range("I122").Value = 3059
$ws.Range("K122").Value = 9177
$ws.Range("M122").Value = -6727

# Row 126: Gold Rush Order
$ws.Range("H126").Value = 9094.429
$ws.Range("I126").Value = 12497.417
$ws.Range("K126").Value = 37492.251
$ws.Range("M126").Value = -35022.251

# Row 132: On Board for Lar
$ws.Range("H132").Value = 775345.9399999999
$ws.Range("I132").Value = 1255875
$ws.Range("K132").Value = 3767625
$ws.Range("M132").Value = -3765095

$ws = $wb.Worksheets.Item("LTW")
# Row 7: Tan Before the Ban
$ws.Range("H7").Value = 50003700
$ws.Range("I7").Value = 250001500
$ws.Range("J7").Value = 4249.875
$ws.Range("K7").Value = 250001500
$ws.Range("L7").Value = 4249.875
$ws.Range("M7").Value = -250001388
$ws.Range("N7").Value = -4473.875

# Row 16: Saddle Sore
$ws.Range("H16").Value = 2201.1428
$ws.Range("I16").Value = 1748.8889
$ws.Range("J16").Value = 3015.2
$ws.Range("K16").Value = 1748.8889
$ws.Range("L16").Value = 3015.2
$ws.Range("M16").Value = -1578.8889
$ws.Range("N16").Value = -3355.2

# Row 40: Best Served Toad
$ws.Range("H40").Value = 3109.8215
$ws.Range("I40").Value = 2707.7273
$ws.Range("K40").Value = 2707.7273
$ws.Range("M40").Value = -2571.7273

# Row 46: Supply Side Logic
$ws.Range("H46").Value = 3307.15
$ws.Range("I46").Value = 1549
$ws.Range("J46").Value = 3502.5
$ws.Range("K46").Value = 1549
$ws.Range("L46").Value = 3502.5
$ws.Range("M46").Value = -1361
$ws.Range("N46").Value = -3878.5

# Row 122: Hell on Leather
$ws.Range("H122").Value = 4669.4443
$ws.Range("I122").Value = 4383.3335
$ws.Range("K122").Value = 13150.0005
$ws.Range("M122").Value = -10700.0005

# Row 126: Battered Books
$ws.Range("H126").Value = 50003700
$ws.Range("I126").Value = 250001500
$ws.Range("J126").Value = 4249.875
$ws.Range("K126").Value = 750004500
$ws.Range("L126").Value = 12749.625
$ws.Range("M126").Value = -750002030
$ws.Range("N126").Value = -17689.625

$ws = $wb.Worksheets.Item("WVR")
# Row 122: Heavy Armoire
$ws.Range("H122").Value = 3316
$ws.Range("I122").Value = 2021.3334
$ws.Range("K122").Value = 6064.0002
$ws.Range("M122").Value = -3614.0002

# Row 126: A Polished Purchase
$ws.Range("H126").Value = 1584.4375
$ws.Range("I126").Value = 1556.7391
$ws.Range("J126").Value = 1655.2222
$ws.Range("K126").Value = 4670.2173
$ws.Range("L126").Value = 4965.6666
$ws.Range("M126").Value = -2200.2173
$ws.Range("N126").Value = -9905.6666
